# Insert a new row for WSTG-AUTH-11 (Testing Multi-Factor Authentication)
# above the current row 49, shifting the "Authorization Testing" section
# (previously starting at row 49) and everything below it down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# Copy the formatting of the row above (row 48, a normal "test" row) only for
# columns A:F, then insert a new row at 49 using that formatting. This avoids
# dragging formatting across the whole 16384-column row, which would otherwise
# inflate the worksheet dimension.
$ws.Range("A48:F48").Copy() | Out-Null
$ws.Range("A49:F49").Insert(-4121, 0) | Out-Null   # xlShiftDown, xlFormatFromLeftOrAbove
$excel.CutCopyMode = 0

# Populate the newly inserted row 49 with the WSTG-AUTH-11 entry.
$ws.Cells.Item(49, 1).Value = $null
$ws.Cells.Item(49, 2).Value = "WSTG-AUTH-11"
$ws.Cells.Item(49, 3).Formula = '=HYPERLINK("https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/04-Authentication_Testing/11-Testing_Multi-Factor_Authentication", "Testing Multi-Factor Authentication (MFA)")'
$ws.Cells.Item(49, 4).Value = "- Identify the type of MFA used by the application.`n- Determine whether the MFA implementation is robust and secure.`n- Attempt to bypass the MFA."
$ws.Cells.Item(49, 5).Value = "Not Started"
$ws.Cells.Item(49, 6).Value = $null

$ws.Rows.Item(49).RowHeight = 66
